$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 440, shifting existing rows 440:499 down to 441:500
$ws.Rows("440:440").Insert()

# Fill in the new row 440 with the new data record
$ws.Cells.Item(440, 1).Value = 1
$ws.Cells.Item(440, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(440, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(440, 4).Value = 45180
$ws.Cells.Item(440, 5).Value = 15
$ws.Cells.Item(440, 6).Value = "Fruta"
$ws.Cells.Item(440, 7).Value = 100102
$ws.Cells.Item(440, 8).Value = "Cítricos"
$ws.Cells.Item(440, 9).Value = 100102003
$ws.Cells.Item(440, 10).Value = "Limón"
$ws.Cells.Item(440, 11).Value = "Tahití"
$ws.Cells.Item(440, 12).Value = "Primera"
$ws.Cells.Item(440, 13).Value = 200
$ws.Cells.Item(440, 14).Value = 50000
$ws.Cells.Item(440, 15).Value = 55000
$ws.Cells.Item(440, 16).Value = 52500
$ws.Cells.Item(440, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(440, 18).Value = "Perú"
$ws.Cells.Item(440, 19).Value = 2188
$ws.Cells.Item(440, 20).Value = 24
